$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 6708.013860684405
$ws.Range("E2").Value = 91228006295.30009
$ws.Range("G2").Value = 91228013005.63039

$ws.Range("B3").Value = 0.2917716402565462
$ws.Range("C3").Value = 0.04071648406533734
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 0.9761466351224579
